# Weekly update: insert a new price record for "Zapallo italiano" at
# Terminal La Palmera de La Serena. The new record is inserted as row 327,
# pushing the existing rows 327-347 down to 328-348 (dimension grows from
# A1:R347 to A1:R348).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 327 (shifts rows 327:347 -> 328:348)
$ws.Rows.Item(327).Insert()

# Populate the new row 327 with the new weekly record
$ws.Range("A327").Value = 8
$ws.Range("B327").Value = "Terminal La Palmera de La Serena"
$ws.Range("C327").Value = "Coquimbo"
$ws.Range("D327").Value = 44746
$ws.Range("E327").Value = 4
$ws.Range("F327").Value = 100112032
$ws.Range("G327").Value = "Zapallo italiano"
$ws.Range("H327").Value = "Sin especificar"
$ws.Range("I327").Value = "Primera"
$ws.Range("J327").Value = 400
$ws.Range("K327").Value = 11000
$ws.Range("L327").Value = 12000
$ws.Range("M327").Value = 11500
$ws.Range("N327").Value = "`$/caja 50 unidades"
$ws.Range("O327").Value = "Región de Arica y Parinacota"
$ws.Range("P327").Value = 230
$ws.Range("Q327").Value = 50
$ws.Range("R327").Value = "Hortaliza"
